$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17: fill in Start time / End time, which drives the Time worked formula,
# and add the description of activities text.
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 15
$ws.Range("E17").Value = "Added info button + pages to each class with respective descriptions"

# Move the active selection to E17 (next empty description cell) like in the edited file.
$ws.Range("E17").Select()
